# Auto-generated Excel COM-interop script
# Applies cached market-data value updates to ALC, ARM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 4148.273
$ws.Range("H111").Value = 1213.6
$ws.Range("H115").Value = 1098
$ws.Range("H127").Value = 901.35297
$ws.Range("H33").Value = 127.111115
$ws.Range("I100").Value = 2997.5
$ws.Range("I111").Value = 756
$ws.Range("I115").Value = 761.6667
$ws.Range("I127").Value = 777.38464
$ws.Range("I33").Value = 116.5
$ws.Range("J100").Value = 4404
$ws.Range("J111").Value = 1900
$ws.Range("J127").Value = 1304.25
$ws.Range("J33").Value = 148.33333
$ws.Range("K100").Value = 2997.5
$ws.Range("K111").Value = 2268
$ws.Range("K115").Value = 2285.0001
$ws.Range("K127").Value = 2332.15392
$ws.Range("K33").Value = 116.5
$ws.Range("L100").Value = 4404
$ws.Range("L111").Value = 5700
$ws.Range("L127").Value = 3912.75
$ws.Range("L33").Value = 148.33333
$ws.Range("M100").Value = -2456.5
$ws.Range("M111").Value = 799
$ws.Range("M115").Value = -718.0001000000002
$ws.Range("M127").Value = 2627.84608
$ws.Range("M33").Value = 112.5
$ws.Range("N100").Value = -5486
$ws.Range("N127").Value = -13832.75
$ws.Range("N33").Value = -606.3333299999999
$ws.Range("N111").Value = -11834

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 31326.666
$ws.Range("H79").Value = 31326.666
$ws.Range("H80").Value = 40960
$ws.Range("H82").Value = 0
$ws.Range("H83").Value = 40960
$ws.Range("H85").Value = 0
$ws.Range("I80").Value = 29800
$ws.Range("I82").Value = 0
$ws.Range("I83").Value = 29800
$ws.Range("I85").Value = 0
$ws.Range("J76").Value = 31326.666
$ws.Range("J79").Value = 31326.666
$ws.Range("J80").Value = 43750
$ws.Range("J82").Value = 0
$ws.Range("J83").Value = 43750
$ws.Range("J85").Value = 0
$ws.Range("K80").Value = 29800
$ws.Range("K82").Value = 0
$ws.Range("K83").Value = 89400
$ws.Range("K85").Value = 0
$ws.Range("L76").Value = 31326.666
$ws.Range("L79").Value = 31326.666
$ws.Range("L80").Value = 43750
$ws.Range("L82").Value = 0
$ws.Range("L83").Value = 131250
$ws.Range("L85").Value = 0
$ws.Range("N76").Value = -32002.666
$ws.Range("N79").Value = -33666.666
$ws.Range("N80").Value = -45746
$ws.Range("N83").Value = -141234
$ws.Range("M80").Value = -28802
$ws.Range("M83").Value = -84408
$ws.Range("M82").ClearContents()
$ws.Range("M85").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 24000
$ws.Range("H31").Value = 17143.932
$ws.Range("H34").Value = 17143.932
$ws.Range("H87").Value = 21999.572
$ws.Range("H90").Value = 21999.572
$ws.Range("I31").Value = 971.2941
$ws.Range("I34").Value = 971.2941
$ws.Range("J100").Value = 24000
$ws.Range("J31").Value = 21967.352
$ws.Range("J34").Value = 21967.352
$ws.Range("J87").Value = 21999.572
$ws.Range("J90").Value = 21999.572
$ws.Range("K31").Value = 971.2941
$ws.Range("K34").Value = 971.2941
$ws.Range("L100").Value = 24000
$ws.Range("L31").Value = 21967.352
$ws.Range("L34").Value = 21967.352
$ws.Range("L87").Value = 21999.572
$ws.Range("L90").Value = 65998.716
$ws.Range("M31").Value = -676.2941
$ws.Range("M34").Value = -769.2941
$ws.Range("N100").Value = -26164
$ws.Range("N31").Value = -22557.352
$ws.Range("N34").Value = -22371.352
$ws.Range("N87").Value = -24371.572
$ws.Range("N90").Value = -77854.716

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2538.8
$ws.Range("I136").Value = 823
$ws.Range("J136").Value = 5112.5
$ws.Range("K136").Value = 2469
$ws.Range("L136").Value = 15337.5
$ws.Range("M136").Value = 2631
$ws.Range("N136").Value = -25537.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 48000
$ws.Range("H15").Value = 23669
$ws.Range("H81").Value = 23669
$ws.Range("H84").Value = 23669
$ws.Range("J130").Value = 48000
$ws.Range("J15").Value = 23669
$ws.Range("J81").Value = 23669
$ws.Range("J84").Value = 23669
$ws.Range("L130").Value = 48000
$ws.Range("L15").Value = 23669
$ws.Range("L81").Value = 23669
$ws.Range("L84").Value = 71007
$ws.Range("N130").Value = -58040
$ws.Range("N15").Value = -24245
$ws.Range("N81").Value = -25665
$ws.Range("N84").Value = -80991

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 47577.332
$ws.Range("H113").Value = 2430
$ws.Range("H114").Value = 31199
$ws.Range("H115").Value = 0
$ws.Range("H116").Value = 29000
$ws.Range("H117").Value = 45000
$ws.Range("H119").Value = 0
$ws.Range("H120").Value = 0
$ws.Range("H130").Value = 34214.5
$ws.Range("H131").Value = 26000
$ws.Range("H55").Value = 222.46666
$ws.Range("H61").Value = 2430
$ws.Range("H82").Value = 1030.7368
$ws.Range("H85").Value = 1030.7368
$ws.Range("I113").Value = 1883.3334
$ws.Range("I131").Value = 26000
$ws.Range("I55").Value = 198.57143
$ws.Range("I61").Value = 1883.3334
$ws.Range("I82").Value = 1033.1666
$ws.Range("I85").Value = 1033.1666
$ws.Range("J110").Value = 47577.332
$ws.Range("J113").Value = 3250
$ws.Range("J114").Value = 31199
$ws.Range("J115").Value = 0
$ws.Range("J116").Value = 29000
$ws.Range("J117").Value = 45000
$ws.Range("J119").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("J130").Value = 34214.5
$ws.Range("J131").Value = 26000
$ws.Range("J55").Value = 243.375
$ws.Range("J61").Value = 3250
$ws.Range("J82").Value = 1029.6154
$ws.Range("J85").Value = 1029.6154
$ws.Range("K113").Value = 1883.3334
$ws.Range("K131").Value = 26000
$ws.Range("K55").Value = 198.57143
$ws.Range("K61").Value = 1883.3334
$ws.Range("K82").Value = 1033.1666
$ws.Range("K85").Value = 1033.1666
$ws.Range("L110").Value = 47577.332
$ws.Range("L113").Value = 3250
$ws.Range("L114").Value = 31199
$ws.Range("L115").Value = 0
$ws.Range("L116").Value = 29000
$ws.Range("L117").Value = 45000
$ws.Range("L119").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("L130").Value = 34214.5
$ws.Range("L131").Value = 26000
$ws.Range("L55").Value = 243.375
$ws.Range("L61").Value = 3250
$ws.Range("L82").Value = 1029.6154
$ws.Range("L85").Value = 1029.6154
$ws.Range("M113").Value = 286.6666
$ws.Range("M55").Value = -25.57142999999999
$ws.Range("M61").Value = -1681.3334
$ws.Range("M82").Value = -672.1666
$ws.Range("M85").Value = 214.8334
$ws.Range("N110").Value = -55757.332
$ws.Range("N113").Value = -7590
$ws.Range("N114").Value = -39877
$ws.Range("N116").Value = -38178
$ws.Range("N117").Value = -54178
$ws.Range("N130").Value = -44254.5
$ws.Range("N131").Value = -36080
$ws.Range("N55").Value = -589.375
$ws.Range("N61").Value = -3654
$ws.Range("N82").Value = -1751.6154
$ws.Range("N85").Value = -3525.6154
$ws.Range("M131").Value = -20960
$ws.Range("N115").ClearContents()
$ws.Range("N119").ClearContents()
$ws.Range("N120").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 446.7143
$ws.Range("H130").Value = 32184.715
$ws.Range("H131").Value = 42049.6
$ws.Range("I107").Value = 425.4
$ws.Range("J107").Value = 500
$ws.Range("J130").Value = 32184.715
$ws.Range("J131").Value = 42049.6
$ws.Range("K107").Value = 1276.2
$ws.Range("L107").Value = 1500
$ws.Range("L130").Value = 32184.715
$ws.Range("L131").Value = 42049.6
$ws.Range("M107").Value = 643.8000000000002
$ws.Range("N107").Value = -5340
$ws.Range("N130").Value = -42224.715
$ws.Range("N131").Value = -52129.6
